# Updated cryptos list on Sat Oct 26 20:30:06 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds text values (e.g. "1.00", "25.60") rather than
# numbers, so force the whole column to Text format before writing any
# values. This prevents Excel from "helpfully" re-interpreting values
# like "1.00" or "25.60" as numbers (which would drop the trailing zero).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.338.61"
$ws.Range("E2").Value = "  +0.73%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.504.01"
$ws.Range("E3").Value = "  +1.04%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "586.64"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6 - Solana
$ws.Range("D6").Value = "172.76"
$ws.Range("E6").Value = "  +2.96%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.502.83"
$ws.Range("E9").Value = "  +0.95%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.27%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.03%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +0.18%  "

# Row 13 - Cardano
$ws.Range("D13").Value = "0.334"
$ws.Range("E13").Value = "  -0.38%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "25.60"
$ws.Range("E14").Value = "  -0.95%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.919.63"

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "66.952.05"
$ws.Range("E16").Value = "  +0.44%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -1.03%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.435.26"
$ws.Range("E18").Value = "  -1.72%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "11.08"
$ws.Range("E19").Value = "  -4.57%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "7.47"
$ws.Range("E20").Value = "  -4.98%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "351.85"
$ws.Range("E21").Value = "  -2.73%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  -0.51%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.06%  "

# Row 24 - NEARProtocol
$ws.Range("D24").Value = "4.26"
$ws.Range("E24").Value = "  -4.23%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "68.69"
$ws.Range("E25").Value = "  -3.03%  "

# Row 26 - SuiNetwork
$ws.Range("E26").Value = "  -1.46%  "

# Row 27 - Aptos
$ws.Range("D27").Value = "9.30"
$ws.Range("E27").Value = "  -1.51%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.08%  "

# Row 29 - WrappedeETH
$ws.Range("D29").Value = "2.619.75"
$ws.Range("E29").Value = "  +0.64%  "

# Row 30 - PEPE
$ws.Range("D30").Value = "0.0₃0911"
$ws.Range("E30").Value = "  -1.98%  "

# Row 31 - Bittensor
$ws.Range("D31").Value = "513.97"
$ws.Range("E31").Value = "  -0.31%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "7.85"
$ws.Range("E32").Value = "  -2.65%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  -2.17%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  -2.81%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.03%  "

# Row 36 & 37 swap: Monero <-> Kaspa
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  -6.50%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "159.83"
$ws.Range("E37").Value = "  +0.98%  "

# Row 38 - WhiteBITCoin
$ws.Range("D38").Value = "18.72"
$ws.Range("E38").Value = "  +0.89%  "

# Row 39 - EthereumClassic
$ws.Range("D39").Value = "18.30"
$ws.Range("E39").Value = "  -3.11%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  -5.03%  "

# Row 41, 42, 43 rotate: USDe/Stacks/RenderToken -> Stacks/RenderToken/USDe
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.69"
$ws.Range("E41").Value = "  -2.84%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "4.86"
$ws.Range("E42").Value = "  -1.80%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.15%  "

# Row 44 - PolygonEcosystemToken
$ws.Range("D44").Value = "0.330"
$ws.Range("E44").Value = "  -1.26%  "

# Row 45 - dogwifhat
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").Value = "  -2.46%  "

# Row 46 - OKB
$ws.Range("D46").Value = "38.80"
$ws.Range("E46").Value = "  -1.08%  "

# Row 47 - Aave
$ws.Range("D47").Value = "143.63"
$ws.Range("E47").Value = "  +0.77%  "

# Row 48 - ARBITRUM
$ws.Range("D48").Value = "0.519"
$ws.Range("E48").Value = "  -3.44%  "

# Row 49 - Filecoin
$ws.Range("D49").Value = "3.47"
$ws.Range("E49").Value = "  -3.18%  "

# Row 50 - BabyDogeCoin
$ws.Range("D50").Value = "0.0₆0253"
$ws.Range("E50").Value = "  -5.46%  "

# Row 51 - Optimism
$ws.Range("D51").Value = "1.58"
$ws.Range("E51").Value = "  -3.92%  "
